$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Station 2 best results")

# Mapping of row number -> new value for column B ("Run 0")
$values = @{
    2 = 0.883116993719734
    3 = 0.7734235047041846
    4 = 0.2104966104216915
    5 = -0.3221380399837128
    6 = 0.2148652162430487
    7 = -0.5554962313654449
    8 = -0.2621492330493018
    9 = 0.4593179363640181
    10 = 0.9005243144425917
    11 = 0.5375381183219355
    12 = -0.8652729942763692
    13 = -0.5630849236071931
    14 = -0.5621503691685834
    15 = -0.4980132539309302
    16 = 0.9420303548197518
    17 = 1.186566923077074
    18 = 1.123697665564034
    19 = -1.173076918193693
    20 = -1.303486083347949
    21 = -1.256745610315601
    22 = -0.39785720593528
    23 = -1.325091554007689
    24 = -1.662810627492533
    25 = -1.160113902378275
    26 = -1.140054108467316
    27 = -0.201913014237361
    28 = -1.256143235221749
    29 = -0.6773372599923947
    30 = 0.3384310006781007
    31 = 1.12419967638983
    32 = 0.8546923913964547
    33 = -0.3276034253762765
    34 = -0.7250406286363725
    35 = -0.4933098132665787
    36 = -0.2430505405447085
    37 = -0.8099329507926127
    38 = -0.9554916999350824
    39 = -0.9667637181708052
    40 = 0.3756838758584351
    41 = 0.3833771199839494
    42 = 0.7223889339836348
    43 = -0.2697466597689135
    44 = 0.513576703916601
    45 = 1.593604802376164
    46 = 2.082910585105029
    47 = 2.344142952628738
    48 = 0.8391724023389752
    49 = -0.1256016537727147
    50 = -0.9552242358195805
    51 = -0.8602662515574869
    52 = -0.7214979173367667
    53 = -1.002630547463725
    54 = -0.2563656776494083
    55 = 0.873120216754867
    56 = 1.617966994459619
    57 = 1.827141363934004
    58 = 1.344934166294925
    59 = 0.7287225506788433
    60 = -0.4938822278882555
    61 = 0.1753018815034183
    62 = 0.6880228715996014
    63 = 0.8982462437409268
    64 = -0.1276902166430932
    65 = 1.869265341871934
    66 = 1.948516875934097
    67 = 2.155132181184887
    68 = 0.4700240045722358
    69 = 0.3902513940204743
    70 = 0.01277845917513423
    71 = 0.09418158463854573
    72 = 0.7150162186583376
    73 = 0.6228138995475769
    74 = 0.6461288016344213
    75 = 0.3383545128598551
    76 = -0.1719166477994486
    77 = 0.4368812041204783
    78 = -0.4594227171984726
    79 = -0.01989217783657342
    80 = 0.6246860855902633
    81 = 1.328402340583235
    82 = 1.234177039130899
    83 = -0.3946763292215997
    84 = -0.2048352273597771
    85 = -0.9145457146101897
    86 = -0.3824600181046837
    87 = -1.488985972521054
    88 = -1.278305294991149
    89 = -0.7051925130673919
    90 = -0.5928362671015492
    91 = 0.3625266146260532
    92 = 0.6671095278253153
    93 = 0.6969174724530497
    94 = 0.4422666916599874
    95 = 0.2736861666587
    96 = -0.3574357399284184
    97 = -1.123450485946833
    98 = -0.7624015581681685
    99 = -0.1766085111652667
    100 = 0.05877926096789207
    101 = -0.1892864444973865
    102 = -0.4587441546929567
    103 = -1.046591949290492
    104 = -1.615465158265159
    105 = -1.841441224158154
    106 = -0.4424715763011001
    107 = -0.5290842981174237
    108 = 0.7229102961871403
    109 = 1.342379514704446
    110 = 1.604159367641256
    111 = 1.100150531102012
    112 = 1.60054010772117
    113 = 1.803769854413962
    114 = 1.695367787105069
    115 = 0.5478833810287235
    116 = 0.4519560278809447
    117 = 1.176694732767472
    118 = 1.554769829227329
    119 = 1.352796755272946
    120 = 0.5299821120704974
    121 = 1.466749946513026
    122 = 1.538773792629724
    123 = 1.921726376658394
    124 = 0.4724825796444003
    125 = 0.9251899844809377
    126 = 0.008268019610376875
    127 = -0.9104032775578831
    128 = 0.1461913671547588
    129 = 0.587273189865653
    130 = 1.818057989194122
    131 = 0.8310994681881189
    132 = 1.171967814701027
    133 = 0.1229289284271944
    134 = 0.117975651568581
    135 = 0.2619056444133643
    136 = 1.132052606855888
    137 = 0.6347000655479004
    138 = -0.2039228855206145
    139 = -0.4857484074710422
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $values[$row]
}
